$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") for rows 2 through 74 was bumped from serial date
# 45181 (2023-09-12) to 45182 (2023-09-13) for every data row.
for ($r = 2; $r -le 74; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
